$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7287194209349384
$ws.Range("C2").Value = 9.226618575922256
$ws.Range("D2").Value = 16.98373111632243
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 273.9243198072813
